$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '303.05'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '5.26%'

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '31.75'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '7.64%'

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.209'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '2.72%'

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '8.93%'

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '7.846'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '6.45%'

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.736'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '8.37%'

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.507'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '7.58%'

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.9082'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-0.85%'

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.01671'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '2,486.68%'

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1682'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '5.22%'

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07487'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '7.90%'

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07947'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '3.25%'

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.02970'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '1.60%'

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.09910'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '10.22%'

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.001487'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-5.38%'

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.04523'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '0.76%'

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.006158'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-0.09%'

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '0.57%'

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '0.02%'

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.1324'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '1.29%'

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.529'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '11.03%'

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.1619'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '2.18%'

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.001217'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '1.80%'

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.004421'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '6.97%'

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '8.25%'

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0001738'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '7.37%'

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04488'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '5.44%'

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007199'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '5.18%'

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1343'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '8.26%'

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '4.41%'

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.01280'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-1.09%'

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00006060'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '5.64%'

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.7068'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-64.04%'

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.01298'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-13.92%'
